# "Se acepta version del servidor en conflictos"
# Rename the report's header columns (row 2) to match the server's field
# names, and drop the stale sample/preview data row (row 3) that is no
# longer part of the report template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: header labels -> server field names
$ws.Range("A2").Value2 = "NroRuc"
$ws.Range("B2").Value2 = "RazonSocial"
$ws.Range("C2").Value2 = "NombreComercial"
$ws.Range("D2").Value2 = "FlagRetencion"
$ws.Range("E2").Value2 = "FlagBloqueoCredito"
$ws.Range("F2").Value2 = "Distrito"
$ws.Range("G2").Value2 = "Direccion"

# Remove the leftover sample data row entirely (shifts dimension to A1:G2)
$ws.Rows("3:3").Delete()
